$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A44").Value = 100
$ws.Range("B44").Value = 2
$ws.Range("C44").Value = 0
$ws.Range("D44").Value = 2
$ws.Range("E44").Value = 0
$ws.Range("F44").Value = 0
$ws.Range("G44").Value = 1
